# Continuação da aula de Validação de Dados: preenche a tabela de exemplo,
# adiciona duas novas linhas (16 e 17), ajusta as validações de dados para
# cobrir o intervalo maior, cria a validação de lista da coluna "Equipe"
# e converte o intervalo A1:F17 em uma Tabela do Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aula 1")
$ws.Activate()

$cpfFormat = '000"."000"."000"-"00'

# --- Linha 2: completa os dados de "Leila Alecrim" ------------------------
$ws.Range("B2").Value = "Mega"
$ws.Range("D2").Value = 53
$ws.Range("E2").Value = 25624

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "05123256984"

# --- Linha 3: CPF de "Filipe Dourado" --------------------------------------
$ws.Range("F3").Value = 12345678963

# Aplica a máscara de CPF em toda a coluna da tabela (linhas 2 a 17)
$ws.Range("F2:F17").NumberFormat = $cpfFormat

# --- Novas linhas 16 e 17 ---------------------------------------------------
$ws.Range("A16").Value = "Junior Caique"
$ws.Range("B17").Value = "Mega"

# --- Ajusta as validações de dados existentes para o novo intervalo --------
$ws.Range("C2:C15").Validation.Delete()
$vC = $ws.Range("C2:C17").Validation
$vC.Add(3, 1, 1, '"Masculino,Feminino,Outros"')
$vC.ErrorTitle = "Informação inválida"
$vC.ErrorMessage = "Você precisa informar um dos seguintes valores: masculino, feminino ou outros."
$vC.InputTitle = "ATENÇÃO"
$vC.InputMessage = "Digite uma das opções: masculino, feminino ou outros."
$vC.ShowInput = $true
$vC.ShowError = $true

$ws.Range("D2:D15").Validation.Delete()
$vD = $ws.Range("D2:D17").Validation
$vD.Add(1, 1, 1, "0", "120")
$vD.ShowInput = $true
$vD.ShowError = $true

$ws.Range("E2:E15").Validation.Delete()
$vE = $ws.Range("E2:E15").Validation
$vE.Add(4, 1, 6, "42369")
$vE.ShowInput = $true
$vE.ShowError = $true

$ws.Range("F2:F15").Validation.Delete()
$vF = $ws.Range("F2:F17").Validation
$vF.Add(6, 1, 3, "11")
$vF.ShowInput = $true
$vF.ShowError = $true

# --- Nova validação de lista para a coluna "Equipe" (B) ---------------------
$vB = $ws.Range("B2:B17").Validation
$vB.Add(3, 1, 1, '$H$2:$H$5')
$vB.ShowInput = $true
$vB.ShowError = $true

# --- Converte o intervalo em Tabela do Excel --------------------------------
$tableRange = $ws.Range("A1:F17")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Tabela1"
$lo.TableStyle = "TableStyleLight8"

# --- Seleção final ----------------------------------------------------------
$ws.Range("C17").Select()
